# Re-split / re-merge w:r run boundaries inside a few paragraphs of the
# document, without altering any visible text or formatting, to match the
# target OOXML structure.
#
# Two primitives are used, both of which operate on Range objects so they
# work purely through the Word COM object model:
#
#   Split-At     -- splits whatever run currently covers $pos into two runs
#                    by toggling a character property (Bold) on/off across
#                    the tail of the paragraph starting at $pos. Word's
#                    engine always materializes a fresh run for a range
#                    whose formatting is set explicitly, even if the
#                    resulting value equals the inherited one, so toggling
#                    Bold true->false leaves the text/formatting unchanged
#                    but forces the run split to stick.
#
#   Collapse-Range -- merges every run inside [start,end] into a single run
#                    by briefly changing the text (appending a sentinel
#                    character) and then restoring the original text. Word
#                    always rebuilds run(s) touched by a real text
#                    replacement, which is exactly what we want to erase
#                    pre-existing run boundaries before re-splitting them
#                    at the new positions.

$d = $word.ActiveDocument

function Split-At($pos, $parEnd) {
    if ($pos -le $parEnd) {
        $tail = $d.Range($pos, $parEnd)
        $tail.Bold = 1
        $tail.Bold = 0
    }
}

function Collapse-Range($startPos, $endPos) {
    $rng = $d.Range($startPos, $endPos)
    $txt = $rng.Text
    $rng.Text = $txt + [char]1
    $rng2 = $d.Range($startPos, $startPos + $txt.Length + 1)
    $rng2.Text = $txt
}

function Resplit-Paragraph($paraIndex, $offsets) {
    $pRange = $d.Paragraphs.Item($paraIndex).Range
    $pStart = $pRange.Start
    $pEnd = $pRange.End - 1
    Collapse-Range $pStart $pEnd
    foreach ($off in $offsets) {
        Split-At ($pStart + $off) $pEnd
    }
}

# --- Paragraph 4: "Our purpose of this descriptive analysis ... station." --
# Split the single run into two at "...their customer |experience...".
$p4 = $d.Paragraphs.Item(4).Range
$p4Start = $p4.Start
$p4End = $p4.End - 1
Split-At ($p4Start + 145) $p4End

# --- Paragraph 6: "That will help us know when ... improve metro rides." --
# Re-split into 7 runs:
#   "That will help us know when to schedule time for cleaning and maintenance. "
#   "Also,"
#   " it will help to bring more people by increasing metro turnstile and"
#   " introduce new technology (digital access)"
#   " making metro rides more convenient, ... most crowded stations and areas"
#   ", and we can analyze riders feedback"
#   " that will help improve metro rides."
Resplit-Paragraph 6 @(75, 80, 148, 190, 334, 370)

# --- Paragraph 10: "We plan to use MTA Turnstile ... 25/09/2021." --
# Re-split into 3 runs:
#   "We plan to use MTA Turnstile "
#   "Data and"
#   " will obtain the data from the MTA.info website. The individual sample ..."
Resplit-Paragraph 10 @(29, 37)

"done"
